# Update the 15 lattice-multiplication problems in the 5x3 table with a new
# set of problems/answers, matching the target OOXML exactly (including the
# xml:space="preserve" attribute on the lines that have leading/trailing
# spaces, and omitting it where there is none - exactly like the original
# document).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each row: table-row, table-col, "AB x CD" problem text,
#           d1,d2 = digits of the second factor (line2 "  d1    d2")
#           e1,e2 = digits of the first factor, reversed (rows "e1|    |" / "e2|    |")
$data = @(
    @(1, 1, "43 x 34", "3", "4", "4", "3"),
    @(1, 2, "85 x 53", "5", "3", "8", "5"),
    @(1, 3, "47 x 13", "1", "3", "4", "7"),
    @(2, 1, "18 x 40", "4", "0", "1", "8"),
    @(2, 2, "27 x 75", "7", "5", "2", "7"),
    @(2, 3, "10 x 39", "3", "9", "1", "0"),
    @(3, 1, "34 x 68", "6", "8", "3", "4"),
    @(3, 2, "68 x 94", "9", "4", "6", "8"),
    @(3, 3, "65 x 48", "4", "8", "6", "5"),
    @(4, 1, "19 x 13", "1", "3", "1", "9"),
    @(4, 2, "84 x 36", "3", "6", "8", "4"),
    @(4, 3, "39 x 43", "4", "3", "3", "9"),
    @(5, 1, "20 x 91", "9", "1", "2", "0"),
    @(5, 2, "99 x 18", "1", "8", "9", "9"),
    @(5, 3, "84 x 30", "3", "0", "8", "4")
)

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr>'
$pkgFooter = '</w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

foreach ($row in $data) {
    $r = $row[0]
    $c = $row[1]
    $problem = $row[2]
    $d1 = $row[3]
    $d2 = $row[4]
    $e1 = $row[5]
    $e2 = $row[6]

    $line2 = "  {0}    {1}" -f $d1, $d2
    $line4 = "{0}|    |" -f $e1
    $line5 = "{0}|    |" -f $e2

    $body = "<w:t>{0}</w:t><w:br/><w:t xml:space=`"preserve`">{1}</w:t><w:br/><w:t xml:space=`"preserve`">  ----</w:t><w:br/><w:t>{2}</w:t><w:br/><w:t>{3}</w:t>" -f $problem, $line2, $line4, $line5

    $xml = $pkgHeader + $body + $pkgFooter

    $cell = $t.Cell($r, $c)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    if ($rng.Start -eq 0) {
        # InsertXML replaces in-place correctly when the range begins at the
        # very start of the document; deleting first here would leave a
        # spurious empty paragraph behind.
        $rng.InsertXML($xml)
    } else {
        # Elsewhere, InsertXML only inserts before Range.End instead of
        # replacing Start..End, so collapse the range via Delete() first.
        $rng.Delete()
        $rng.InsertXML($xml)
    }
}
